$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" — the handoff xliff generation timestamp for
# e2e\6fa58431-d3a6-4c97-a1f5-b4c604f0972c.md was refreshed to a new value.
$newDate = "2016-08-29 00:41:14"

# Overview sheet: row 6 is the 6fa58431-... file; column G is
# "Latest HO Xliff Generate Date".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value2 = $newDate

# de-de sheet: rows 5 (2813a0d6-...) and 6 (6fa58431-...) share the same
# "Latest Handoff Datetime" (column H) batch timestamp and both move forward.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value2 = $newDate
$wsDeDe.Range("H6").Value2 = $newDate
